# Rebuild this inventory-upload template for a new product ("拉米夫定片"
# instead of the old "促肝细胞生长素肠溶胶囊" sample row), and rename the
# sheet to the generic "Sheet1".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sheet1"

# Header row: last column header becomes "库存数量" (quantity in stock).
$ws.Range("D1").Value = "库存数量"

# Sample/example data row. Item code must stay text so the leading zero in
# "062960" survives (otherwise Excel would coerce it to the number 62960).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "062960"
$ws.Range("B2").Value = "拉米夫定片"
$ws.Range("C2").Value = "100mgX14s"

# Drop the thin-box border that used to outline this sample row (D2, with
# its numeric 100, keeps its original bordered style).
$ws.Range("A2:C2").Borders.LineStyle = -4142

# Widen the columns for the new, longer sample content.
$ws.Columns.Item(1).ColumnWidth = 14
$ws.Columns.Item(2).ColumnWidth = 33.285714285714285
$ws.Columns.Item(3).ColumnWidth = 19.142857142857142
$ws.Columns.Item(4).ColumnWidth = 24.285714285714285

$ws.Cells.Item(10, 2).Select() | Out-Null
